# Update cryptocurrency price (D) and volume change (E) columns
# to reflect the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.874.68"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.643.35"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "218.59"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.870.32"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "1.649.77"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "65.51"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "26.872.03"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").Value = "215.28"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  +6.71%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "2.37"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Value = "147.65"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "0.119"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "15.79"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "1.283.46"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "0.534"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "0.817"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "5.36"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "1.782.50"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("D45").Value = "92.69"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").Value = "61.17"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "7.59"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").Value = "  -0.13%  "
